$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Purchase 22-23": a new ledger line is inserted under Sr. No 1's
# block (after row 8) recording a further deduction of 100,000 from the
# running total. Every row from the old row 9 downward shifts down by one.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a blank row at position 9 (pushes old rows 9:25 down to 10:26).
$ws1.Rows("9:9").Insert()

# The new row should look like the rest of its section (row 8's styling:
# regular font for A:E, thin borders, etc.) so copy formatting down first.
$ws1.Range("A8:F8").Copy()
$ws1.Range("A9:F9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 8's cell F8 was previously the bold "block total" cell; now that the
# block continues into row 9, it becomes a plain running-total cell like
# F6/F7 above it.
$ws1.Range("F7").Copy()
$ws1.Range("F8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New row 9 is the block's closing/total line: prior running total less a
# 100,000 adjustment/write-off.
$ws1.Range("F9").Formula = "=F8-100000"

# Leave the selection where the user last left it.
[void]$ws1.Range("I16").Select()

# ---------------------------------------------------------------------------
# Sheet "Sale 22-23": the header row (row 29) had an oversized manual row
# height; reset it back to the sheet's normal auto-calculated height.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows(29).AutoFit()
